# Addressable 적용 작업 1차
#
# The ItemToolTipData sheet's "iconAssetPath" column (B) held bare asset
# paths (e.g. "Sprites/ItemIcon/Prop_A"); switch them to explicit file
# names with the ".png" extension so the Addressable loader can resolve
# them. Every other cell in the workbook that referenced the shared
# string table shifts automatically as Excel repacks sst.xml - no other
# sheet needs touching for data.
#
# Cells are updated in the same order the original author's Excel
# session produced them (bottom groups first) so the appended shared
# -string order in the saved package lines up with the source edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ItemToolTipData")

$rows = @(22, 23, 20, 21, 19, 18, 17, 16, 15, 14, 13, 12, 11, 10, 9, 3, 4, 5, 6, 7, 8)
foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 2)
    $old = $cell.Value2
    $cell.Value = $old + ".png"
}

# The author's focus moved from the PropItemData tab to ItemToolTipData
# (workbook.xml activeTab 2 -> 3), landing the selection on C11.
$ws.Activate() | Out-Null
$ws.Range("C11").Select() | Out-Null
